# Results including qb and painter.
# Extends the "TROY" sheet with two additional benchmark blocks (painter
# tasks/node=2 at rows 16-19, qb tasks/node=4 at rows 25-28), tweaks the
# existing table's header/TTC wording and the BJ-Diane throughput number,
# and records a couple of free-form annotations (error note, timings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 1. Tweaks to the existing table (rows 6-9)
# ---------------------------------------------------------------------
# Header label becomes more specific.
$ws.Range("H6").Value = "TTC(only matching)"

# Measured throughput for BJ-Diane was re-run -> updated number, with a
# note about the run that needed 4 workers/node.
$ws.Range("H8").Value = 941
$ws.Range("K8").Value = "error with four workers per node"

# ---------------------------------------------------------------------
# 2. New block: painter, 2 tasks/node (rows 16-19)
# ---------------------------------------------------------------------
$ws.Range("B6:H6").Copy()
$ws.Range("B16:H16").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 58

$ws.Range("B16").Value = "Backend"
$ws.Range("C16").Value = "Number of tasks"
$ws.Range("D16").Value = "Number of cores"
$ws.Range("E16").Value = "Number of threads per task"
$ws.Range("F16").Value = "Machine"
$ws.Range("G16").Value = "Type"
$ws.Range("H16").Value = "TTC(only matching)"

$ws.Range("K16").Value = "tasks/node"
$ws.Range("K16").Font.Bold = $true
$ws.Range("K16").Font.Size = 15
$ws.Range("K16").WrapText = $true

$ws.Rows.Item(17).RowHeight = 16

$ws.Range("B17").Value = "BJ-SAGA"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = "painter"
$ws.Range("G17").Value = "matching"
$ws.Range("H17").Value = 529
$ws.Range("K17").Value = 2

$ws.Range("B18").Value = "BJ-Diane"
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 16
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = "painter"
$ws.Range("G18").Value = "matching"
$ws.Range("H18").Value = 542
$ws.Range("K18").Value = 2

$ws.Range("B19").Value = "BJ-SAGA, BJ-Diane"
$ws.Range("C19").Value = "4,4"
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = "painter"
$ws.Range("G19").Value = "matching"
$ws.Range("H19").Value = 545
$ws.Range("K19").Value = 2

# ---------------------------------------------------------------------
# 3. New block: qb, 4 tasks/node (rows 25-28)
# ---------------------------------------------------------------------
$ws.Range("B6:H6").Copy()
$ws.Range("B25:H25").PasteSpecial(-4122)
$ws.Rows.Item(25).RowHeight = 58

$ws.Range("B25").Value = "Backend"
$ws.Range("C25").Value = "Number of tasks"
$ws.Range("D25").Value = "Number of cores"
$ws.Range("E25").Value = "Number of threads per task"
$ws.Range("F25").Value = "Machine"
$ws.Range("G25").Value = "Type"
$ws.Range("H25").Value = "TTC(only matching)"

$ws.Range("J25").Value = "real time"
$ws.Range("J25").Font.Bold = $true
$ws.Range("J25").Font.Size = 15
$ws.Range("J25").WrapText = $true

$ws.Range("K25").Value = "tasks/node"
$ws.Range("K25").Font.Bold = $true
$ws.Range("K25").Font.Size = 15
$ws.Range("K25").WrapText = $true

$ws.Rows.Item(26).RowHeight = 16

$ws.Range("B26").Value = "BJ-SAGA"
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = "qb"
$ws.Range("G26").Value = "matching"
$ws.Range("H26").Value = 494
$ws.Range("J26").Value = "8m17.946s"
$ws.Range("K26").Value = 4

$ws.Range("B27").Value = "BJ-Diane"
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 16
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = "qb"
$ws.Range("G27").Value = "matching"
$ws.Range("H27").Value = 485
$ws.Range("J27").Value = "12m11.266s`n"
$ws.Range("J27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 45
$ws.Range("K27").Value = 4

$ws.Range("B28").Value = "BJ-SAGA, BJ-Diane"
$ws.Range("C28").Value = "4, 4"
$ws.Range("D28").Value = 16
$ws.Range("E28").Value = 2
$ws.Range("F28").Value = "qb"
$ws.Range("G28").Value = "matching"
$ws.Range("H28").Value = 521
$ws.Range("J28").Value = "10m17.069s"
$ws.Range("K28").Value = 4

# ---------------------------------------------------------------------
# 4. Misc view tweaks
# ---------------------------------------------------------------------
$ws.Range("H31").Select()
